$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 374 (dimension A1:D374).
# We need to append 11 new rows (375:385) continuing the daily series,
# matching the formatting of the last existing data row (374).

# Copy formatting (style) of the last populated row down into the new rows
$ws.Range("A374:D374").Copy()
$ws.Range("A375:D385").PasteSpecial(-4122)  # xlPasteFormats

$dates = @(44449,44450,44451,44452,44453,44454,44455,44456,44457,44458,44459)
$bVals = @(0,0,0,0,0,0,0,0,1,0,0)
$cVals = @(0,0,0,0,0,0,0,0,1,1,1)
$dVals = @(0,0,0,0,0,0,0,0,62.34413965087282,62.34413965087282,62.34413965087282)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 375 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}
